$wb = $excel.ActiveWorkbook

$wsSSSOM = $wb.Worksheets.Item("SSSOM")

# Update existing mapping row's predicate from skos:closematch to skos:narrowmatch
$wsSSSOM.Range("C2").Value = "skos:narrowmatch"

# Add a new mapping row (row 3) for FM1.3
$wsSSSOM.Range("A3").Value = "estuarine:Estuarine, delta"
$wsSSSOM.Range("B3").Value = "Estuarine, delta"
$wsSSSOM.Range("C3").Value = "skos:narrowmatch"
$wsSSSOM.Range("D3").Value = "get:groups/M1.3"
$wsSSSOM.Range("E3").Value = "FM1.3 Intermittently closed and open lakes and lagoons"
$wsSSSOM.Range("F3").Value = "semapv:ManualMappingCuration"
$wsSSSOM.Range("G3").Value = "orcid:0009-0001-6090-9959"
$wsSSSOM.Range("H3").Value = "Craig Macfarlane"
$wsSSSOM.Range("I3").Value = Get-Date -Year 2024 -Month 5 -Day 10 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$wsSSSOM.Range("K3").Value = "status:draft"

# Make SSSOM the active sheet/tab, with E3 selected
$wsSSSOM.Activate()
$wsSSSOM.Range("E3").Select()
